# work in progress on vaporize_theia.py
# - remove the "KCaAlSi2O7" (R40) reaction row from Table 3 and renumber the
#   trailing reaction labels so they stay contiguous (R40..R46)
# - fix the ZnO_g reactant description on Table 4 (Zn -> Zn_l)

$wb = $excel.ActiveWorkbook

# --- Table 3: drop the KCaAlSi2O7 reaction (row 41) ---------------------
$ws3 = $wb.Worksheets.Item("Table 3")
$ws3.Select()

$ws3.Rows(41).Delete()

for ($i = 41; $i -le 47; $i++) {
    $ws3.Cells.Item($i, 1).Value = "R" + ($i - 1)
}

$ws3.Range("D28").Select()

# --- Table 4: 1*Zn, 0.5*O2 -> 1*Zn_l, 0.5*O2 on the ZnO_g row ------------
$ws4 = $wb.Worksheets.Item("Table 4")
$ws4.Select()

$ws4.Range("E36").Value = "1*Zn_l, 0.5*O2"

$ws4.Range("H8").Select()
